# Workshop schedule update
#  - Replace the generic "Linköping" room entry with specific room names
#    for each day that needed one, and add a "link_room" (column K) popup
#    link to the MazeMap location for that room.
#  - (Delimiter / iterateDelta cleanup happens implicitly on save.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Room names (column B) for the days that have a specific room assigned.
$roomPapaver = "Room Papaver, Hus 511/001, LiU"
$roomDolomit = "Room Dolomit, Hus 440, LiU"
$roomAntracit = "Room Antracit, Hus 440, LiU"

$ws.Range("B2").Value  = $roomPapaver
$ws.Range("B8").Value  = $roomPapaver
$ws.Range("B15").Value = $roomDolomit
$ws.Range("B19").Value = $roomDolomit
$ws.Range("B25").Value = $roomAntracit

# Room map links (column K / link_room) - one MazeMap link per room, added
# as a clickable hyperlink so the text shown is the URL itself.
$ws.Hyperlinks.Add($ws.Range("K2"),  "https://link.mazemap.com/00mnumNU")
$ws.Hyperlinks.Add($ws.Range("K8"),  "https://link.mazemap.com/00mnumNU")
$ws.Hyperlinks.Add($ws.Range("K15"), "https://link.mazemap.com/up3GnjPm")
$ws.Hyperlinks.Add($ws.Range("K19"), "https://link.mazemap.com/up3GnjPm")
$ws.Hyperlinks.Add($ws.Range("K25"), "https://link.mazemap.com/wXeFDYNR")

# Restore the cursor/selection position saved with the workbook.
$ws.Range("B26").Select()
